# "scorpion laten draaien en bugs eruit"
# Log a new "Donderdag" (Thursday) entry in the "week 3" sheet: the
# scorpion was run and bugs were pulled out of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 3")

# Day label (reuses the existing "Donderdag" shared string used elsewhere
# in the log).
$ws.Range("A11").Value = "Donderdag"

# Date column: copy the number formatting from the previous "new day" row
# (B7) so we land on the same date style instead of minting a new one,
# then fill in the date serial (23 Jan 2014).
$ws.Range("B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B11").Value = 41662

# Begin/eind tijd + id for the new activity entry.
$ws.Range("C11").Value = 0.37152777777777773
$ws.Range("D11").Value = 0.3979166666666667
$ws.Range("E11").Value = 5

# Activiteiten description.
$ws.Range("F11").Value = "De scorpion laten draaien en foute eruit gehaald."

# The row grows to fit the wrapped description text.
$ws.Rows.Item(11).RowHeight = 28.5

# Move the active selection to the newly-edited cell.
$ws.Range("F11").Select() | Out-Null
